# "ran scraper to update the data"
#
# The scraper run adds two new stat columns - "height" and "weight" - to
# the sheet, pushing the existing "fantasy points" column two slots to the
# right (from E to G) and replacing it with freshly scraped "height" values
# while the brand-new "weight" column lands in between.
#
# End result:
#   B1:E1 headers shift meaning only for E1 (now "height"); F1 becomes
#   "weight"; G1 becomes "fantasy points" (the old E1 header, relocated).
#   Column E's old per-row "fantasy points" numbers move into column G.
#   Column E is refreshed with the new scraped "height" value (constant
#   across all rows this run: 6.333333333333333).
#   Column F is populated with the new scraped "weight" value (constant
#   across all rows this run: 253).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$headerRow = $firstRow

$oldFantasyCol = 5   # E
$newHeightCol  = 5   # E (reused / overwritten in place)
$newWeightCol  = 6   # F (new)
$movedFantasyCol = 7 # G (new - receives the old E data)

# --- Relocate the existing "fantasy points" data (column E) out to the new
#     column G before E gets overwritten with the refreshed scrape values. ---
for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $oldValue = $ws.Cells.Item($r, $oldFantasyCol).Value2
    $ws.Cells.Item($r, $movedFantasyCol).Value = $oldValue
}

# --- Give the two new header cells (F1, G1) the same look as the other
#     styled header cells by copying an existing header cell's formatting
#     onto them, then set the correct label text afterwards. ---
$ws.Cells.Item($headerRow, $oldFantasyCol).Copy($ws.Range($ws.Cells.Item($headerRow, $newWeightCol), $ws.Cells.Item($headerRow, $movedFantasyCol)))

$ws.Cells.Item($headerRow, $newHeightCol).Value = "height"
$ws.Cells.Item($headerRow, $newWeightCol).Value = "weight"
$ws.Cells.Item($headerRow, $movedFantasyCol).Value = "fantasy points"

# --- Write the freshly scraped height/weight values for every data row. ---
$ws.Range($ws.Cells.Item($headerRow + 1, $newHeightCol), $ws.Cells.Item($lastRow, $newHeightCol)).Value = 6.333333333333333
$ws.Range($ws.Cells.Item($headerRow + 1, $newWeightCol), $ws.Cells.Item($lastRow, $newWeightCol)).Value = 253
